# Apply the 2023-04-13 cryptos-list refresh (Price column D, Volume(1h) column E).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column E ("Volume(1h)") - plain text cells, safe to assign directly ---
$ws.Range("E2").Value = "  +0.76%  "
$ws.Range("E3").Value = "  +5.38%  "
$ws.Range("E4").Value = "  -0.48%  "
$ws.Range("E5").Value = "  +0.88%  "
$ws.Range("E7").Value = "  +1.35%  "
$ws.Range("E8").Value = "  +3.35%  "
$ws.Range("E9").Value = "  +2.62%  "
$ws.Range("E10").Value = "  +3.56%  "
$ws.Range("E11").Value = "  +0.79%  "
$ws.Range("E12").Value = "  +3.14%  "
$ws.Range("E13").Value = "  +2.69%  "
$ws.Range("E14").Value = "  +4.50%  "
$ws.Range("E15").Value = "  +2.63%  "
$ws.Range("E16").Value = "  +0.14%  "
$ws.Range("E17").Value = "  +2.10%  "
$ws.Range("E18").Value = "  +1.39%  "
$ws.Range("E19").Value = "  +1.67%  "
$ws.Range("E20").Value = "  +3.59%  "
$ws.Range("E21").Value = "  -0.28%  "
$ws.Range("E22").Value = "  +3.32%  "
$ws.Range("E23").Value = "  +0.86%  "
$ws.Range("E24").Value = "  +2.93%  "
$ws.Range("E25").Value = "  +1.38%  "
$ws.Range("E26").Value = "  +6.37%  "
$ws.Range("E27").Value = "  +4.34%  "
$ws.Range("E28").Value = "  +1.28%  "
$ws.Range("E29").Value = "  +6.96%  "
$ws.Range("E30").Value = "  +2.14%  "
$ws.Range("E31").Value = "  +6.42%  "
$ws.Range("E32").Value = "  +1.93%  "
$ws.Range("E33").Value = "  +1.64%  "
$ws.Range("E34").Value = "  +4.28%  "
$ws.Range("E35").Value = "  +1.50%  "
$ws.Range("E36").Value = "  +10.47%  "
$ws.Range("E37").Value = "  +2.93%  "
$ws.Range("E38").Value = "  +2.35%  "
$ws.Range("E39").Value = "  +1.42%  "
$ws.Range("E40").Value = "  +5.16%  "
$ws.Range("E41").Value = "  +3.87%  "
$ws.Range("E42").Value = "  +4.22%  "
$ws.Range("E43").Value = "  +0.31%  "
$ws.Range("E44").Value = "  +4.93%  "
$ws.Range("E45").Value = "  +3.60%  "
$ws.Range("E46").Value = "  +4.72%  "
$ws.Range("E47").Value = "  +0.23%  "
$ws.Range("E48").Value = "  +0.68%  "
$ws.Range("E49").Value = "  +1.44%  "
$ws.Range("E50").Value = "  +2.28%  "
$ws.Range("E51").Value = "  +1.12%  "

# --- Column D ("Price") values that are not parseable as a single number
# (thousand-grouped, look like "30.206.79") stay text automatically. ---
$ws.Range("D2").Value = "30.206.79"
$ws.Range("D3").Value = "1.970.47"
$ws.Range("D14").Value = "1.950.30"
$ws.Range("D23").Value = "30.231.80"
$ws.Range("D27").Value = "2.172.91"

# --- Column D ("Price") values that parse as plain decimals would silently
# become numeric cells via a straight .Value assignment, losing the exact
# text the source keeps (e.g. "6.040", "1.004", "0.00001102"). Force each one
# to stay text by temporarily marking the cell as Text before assigning, then
# restore the default "Normal" style so no stray formatting is left behind.
# (Looping per-cell rather than building one multi-area Range, since applying
# a property to a comma-joined Range only touches its first area here.) ---
$priceRows = @(4, 5, 7, 8, 9, 10, 11, 12, 13, 15, 16, 17, 18, 19, 20, 21, 22, 24, 25, 26, 28, 29, 30, 31, 32, 33, 34, 35, 36, 37, 39, 40, 41, 42, 43, 44, 45, 46, 47, 48, 49, 50, 51)
$priceVals = @("0.9979", "322.02", "0.5115", "0.4092", "0.08421", "1.131", "42.46", "24.19", "6.462", "7.375", "1.004", "93.78", "0.00001102", "0.06532", "18.76", "0.9986", "6.040", "11.46", "2.196", "22.63", "162.65", "2.364", "129.92", "1.141", "0.1053", "6.035", "3.789", "0.02474", "1.295", "5.364", "0.2171", "8.918", "0.6551", "11.74", "1.224", "13.52", "0.6117", "2.184", "3.631", "123.41", "1.220", "79.22", "1.129")
for ($i = 0; $i -lt $priceRows.Length; $i++) {
    $cell = $ws.Range("D" + $priceRows[$i])
    $cell.NumberFormat = "@"
    $cell.Value = $priceVals[$i]
    $cell.Style = "Normal"
}

